$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.172.89"
$ws.Range("E2").Value = "  +0.14%  "

$ws.Range("D3").Value = "2.522.96"
$ws.Range("E3").Value = "  +0.30%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.81%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.01%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("E8").Value = "  +0.37%  "

$ws.Range("D9").Value = "2.521.61"
$ws.Range("E9").Value = "  +0.13%  "

$ws.Range("E10").Value = "  +0.82%  "

$ws.Range("E11").Value = "  +0.09%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.35"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.45%  "

$ws.Range("E13").Value = "  -2.98%  "

$ws.Range("D14").Value = "2.948.31"
$ws.Range("E14").Value = "  -0.44%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.19"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.46%  "

$ws.Range("D16").Value = "58.996.80"
$ws.Range("E16").Value = "  -0.06%  "

$ws.Range("E17").Value = "  -0.65%  "

$ws.Range("D18").Value = "2.520.58"
$ws.Range("E18").Value = "  +0.22%  "

$ws.Range("E19").Value = "  +0.74%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.72%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "326.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.20%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.14%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.88"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.75%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.75"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.63%  "

$ws.Range("E25").Value = "  -0.38%  "

$ws.Range("E26").Value = "  +0.18%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.03%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.66"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.91%  "

$ws.Range("D29").Value = "0.0₃0778"
$ws.Range("E29").Value = "  +0.85%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.69"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.01%  "

$ws.Range("E31").Value = "  -0.88%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "169.60"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.91%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.20"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.75%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.77%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.54"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.41%  "

$ws.Range("E37").Value = "  -2.72%  "

$ws.Range("E38").Value = "  -1.13%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.66"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.89%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.824"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.35%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.63"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.33%  "

$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "285.21"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.85%  "

$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.24"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.14%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.12%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "131.14"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.27%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.604"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.73%  "

$ws.Range("E47").Value = "  +0.13%  "

$ws.Range("E48").Value = "  +0.35%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0514"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.86%  "

$ws.Range("E50").Value = "  -0.60%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.56"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.32%  "
